$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 2.3.6: the trailing description is removed, leaving just "2.3.6 ", and a
# second run containing " X" (the "done" marker used throughout this
# checklist document) is appended right after it.
# ---------------------------------------------------------------------------

$target = "2.3.6 thêm chưa hoàn thiện, còn xóa và cập nhật"

$findRange = $d.Content
$found = $findRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Grab the whole paragraph (including its trailing paragraph mark) so the
    # replacement XML below lands back in a single paragraph instead of
    # spawning a new one. Re-materialising the range via $d.Range(start, end)
    # (rather than calling InsertXML straight off the chained
    # Paragraphs(1).Range) keeps the paragraph's own <w:pPr> intact.
    $rawParaRange = $findRange.Paragraphs(1).Range
    $paraRange = $d.Range($rawParaRange.Start, $rawParaRange.End)

    $newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">2.3.6 </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> X</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $paraRange.InsertXML($newParaXml)
}
